# Auto-generated edit script applying the crypto price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks numeric: force Text formatting first so Excel
# does not silently coerce the string into a Number, then restore the
# original (default/"Normal") style so no stray formatting is introduced.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9979"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9994"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3926"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3852"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.381"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9979"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "49.41"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08466"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.00"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.068"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.563"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001282"
$ws.Range("D16").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "93.71"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06920"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.818"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9995"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.45"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.432"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.907"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "156.79"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "139.72"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.262"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.910"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.490"
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08115"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9893"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02894"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.615"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2680"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.09171"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.37"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "13.67"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7530"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.13"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6916"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.481"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.074"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9985"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08260"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "133.88"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.226"
$ws.Range("D51").Style = "Normal"

# Plain text assignments (values that Excel will not misinterpret as numbers)
$ws.Range("D2").Value = "23.886.23"
$ws.Range("E2").Value = "  -3.00%  "
$ws.Range("D3").Value = "1.621.79"
$ws.Range("E3").Value = "  -3.18%  "
$ws.Range("E4").Value = "  -0.43%  "
$ws.Range("E5").Value = "  -1.91%  "
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("E7").Value = "  -0.40%  "
$ws.Range("E8").Value = "  -2.54%  "
$ws.Range("B9").Value = "Polygon"
$ws.Range("C9").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("E9").Value = "  -1.64%  "
$ws.Range("B10").Value = "BinanceUSD"
$ws.Range("C10").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("E10").Value = "  -0.42%  "
$ws.Range("E11").Value = "  -2.80%  "
$ws.Range("E12").Value = "  -2.15%  "
$ws.Range("E13").Value = "  -4.99%  "
$ws.Range("E14").Value = "  -3.76%  "
$ws.Range("E15").Value = "  -2.04%  "
$ws.Range("E16").Value = "  -2.71%  "
$ws.Range("D17").Value = "1.621.40"
$ws.Range("E17").Value = "  -3.63%  "
$ws.Range("E18").Value = "  -0.37%  "
$ws.Range("E19").Value = "  -1.41%  "
$ws.Range("E20").Value = "  -4.80%  "
$ws.Range("E21").Value = "  -3.76%  "
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("E23").Value = "  -3.54%  "
$ws.Range("D24").Value = "23.876.91"
$ws.Range("E24").Value = "  -3.09%  "
$ws.Range("E25").Value = "  +3.44%  "
$ws.Range("E26").Value = "  +4.53%  "
$ws.Range("E27").Value = "  -3.79%  "
$ws.Range("E28").Value = "  -2.04%  "
$ws.Range("E29").Value = "  -4.56%  "
$ws.Range("E30").Value = "  -10.54%  "
$ws.Range("E31").Value = "  -5.60%  "
$ws.Range("E32").Value = "  -0.22%  "
$ws.Range("D33").Value = "1.795.77"
$ws.Range("E33").Value = "  -3.69%  "
$ws.Range("E34").Value = "  -2.37%  "
$ws.Range("E35").Value = "  -0.31%  "
$ws.Range("E36").Value = "  -6.34%  "
$ws.Range("E37").Value = "  -5.11%  "
$ws.Range("E38").Value = "  -4.50%  "
$ws.Range("E39").Value = "  -4.77%  "
$ws.Range("E40").Value = "  +0.78%  "
$ws.Range("E41").Value = "  +0.97%  "
$ws.Range("E42").Value = "  -6.88%  "
$ws.Range("E43").Value = "  -4.95%  "
$ws.Range("E44").Value = "  -3.11%  "
$ws.Range("E45").Value = "  -2.78%  "
$ws.Range("E46").Value = "  -3.24%  "
$ws.Range("E47").Value = "  -2.42%  "
$ws.Range("E48").Value = "  -0.21%  "
$ws.Range("E49").Value = "  -4.60%  "
$ws.Range("E50").Value = "  -2.97%  "
$ws.Range("E51").Value = "  -7.76%  "
